$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("Total")
$wsNonCollab = $wb.Worksheets.Item("Non-Collaborative")

# 1111MCL upgraded to ACO Insight Enhanced - update their monthly fee
# from 1995 to 2150 on the Non-Collaborative pricing sheet.
$wsNonCollab.Range("B2").Formula = "=2500+2150"

for ($r = 3; $r -le 32; $r++) {
    $wsNonCollab.Cells.Item($r, 2).Value = 2150
}

# Leave the cursor where the editor left it on this sheet, then
# re-activate the Total sheet so it stays the active tab.
$wsNonCollab.Range("E11").Select()
$wsTotal.Activate()
